$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-18 from serial date 45190 to 45192
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45192
}
